# Select the "ProductLoanInput" sheet, update the product id value, and make
# this sheet the active one (so the workbook re-opens on it instead of
# "ProductLoanOutput").

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the product id on ProductLoanInput!B2
$wsInput.Range("B2").Value = 3533

# Activate ProductLoanInput and select B2 as the active cell, so the
# workbook re-opens on this sheet (instead of "ProductLoanOutput") with
# B2 selected.
$wsInput.Activate()
$wsInput.Range("B2").Select()
